$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Fgf1", "Nrp1", "ECs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 61.9442646948221, 557.498382253399, 0.02090678169092003, 0.02090678169092003)
    ,@("ECs", "Fgf1", "Nrp1", "FAPs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 31.15178269536577, 280.366044258292, 0.0105140245558429, 0.01051402455584291)
    ,@("ECs", "Fgf1", "Nrp1", "M2", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 29.82901169675678, 268.461105270811, 0.01006757669450748, 0.01006757669450749)
    ,@("ECs", "Fgf1", "Nrp1", "sCs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 10.08200213763445, 90.73801923871001, 0.003402772132938589, 0.003402772132938589)
    ,@("FAPs", "Fgf1", "Nrp1", "ECs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 534.8892708202409, 4814.003437382167, 0.1805302439047104, 0.1805302439047104)
    ,@("FAPs", "Fgf1", "Nrp1", "FAPs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 268.9959177458382, 2420.963259712544, 0.09078869457515754, 0.09078869457515755)
    ,@("FAPs", "Fgf1", "Nrp1", "M2", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 257.5737785309502, 2318.164006778552, 0.08693361336327385, 0.08693361336327385)
    ,@("FAPs", "Fgf1", "Nrp1", "sCs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 87.05817719163556, 783.5235947247199, 0.02938296731621547, 0.02938296731621547)
    ,@("M2", "Fgf1", "Nrp1", "ECs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 9.444036429695444, 84.996327867259, 0.003187452605813928, 0.003187452605813928)
    ,@("M2", "Fgf1", "Nrp1", "FAPs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 4.74940774701911, 42.744669723172, 0.00160297053193357, 0.00160297053193357)
    ,@("M2", "Fgf1", "Nrp1", "M2", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 4.547737785150111, 40.929640066351, 0.001534904991286935, 0.001534904991286935)
    ,@("M2", "Fgf1", "Nrp1", "sCs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 1.537104297567778, 13.83393867811, 0.0005187873993459596, 0.0005187873993459597)
    ,@("sCs", "Fgf1", "Nrp1", "ECs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 773.5992566292183, 6962.393309662964, 0.2610971468349194, 0.2610971468349195)
    ,@("sCs", "Fgf1", "Nrp1", "FAPs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 389.0432157769133, 3501.38894199222, 0.1313058056408135, 0.1313058056408136)
    ,@("sCs", "Fgf1", "Nrp1", "M2", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 372.5236127715984, 3352.712514944385, 0.1257302816539794, 0.1257302816539794)
    ,@("sCs", "Fgf1", "Nrp1", "sCs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 125.9104357349834, 1133.19392161485, 0.04249597610834101, 0.04249597610834102)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}

Write-Output "Updated $($data.Count) data rows"
